# PTW: AXI4 Redesigned, Test added
# Applies the Sheet2 "Fetch" truth-table additions + a couple of
# corrections to the existing rows, as captured by the commit's xlsx diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Fix up existing row 9 ---------------------------------------------
$ws.Range("C9").Value = "NONE"
$ws.Range("G9").Value = 1

# --- New row 10 -----------------------------------------------------------
$ws.Range("A10").Value = "BRANCH_TAKEN"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "NONE"
$ws.Range("D10").Value = "IDLE"
$ws.Range("E10").Value = "INTERRUPT_PENDING"
$ws.Range("F10").Value = "ready"
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0

# --- New row 11 -----------------------------------------------------------
$ws.Range("A11").Value = "BRANCH_TAKEN"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "NONE"
$ws.Range("D11").Value = "IDLE"
$ws.Range("E11").Value = "NONE"
$ws.Range("F11").Value = "abort"
$ws.Range("G11").Value = 1
$ws.Range("I11").Value = "abort because IRQ handling"

# --- New row 12 -----------------------------------------------------------
$ws.Range("A12").Value = "BRANCH_TAKEN"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "MTVEC"
$ws.Range("D12").Value = "IDLE"
$ws.Range("E12").Value = "NONE"
$ws.Range("F12").Value = "branch taken"
$ws.Range("G12").Value = 0

# --- New row 13 -----------------------------------------------------------
$ws.Range("A13").Value = "MTVEC"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = "NONE"
$ws.Range("D13").Value = "DONE"
$ws.Range("E13").Value = "INSTR"
$ws.Range("F13").Value = "not ready"
$ws.Range("G13").Value = 0

# --- New row 14 -----------------------------------------------------------
$ws.Range("D14").Value = "IDLE"
$ws.Range("E14").Value = "INSTR"
$ws.Range("F14").Value = "ready"
$ws.Range("G14").Value = 0

# --- Column widths (E and I got wider / no longer auto "best fit") --------
$ws.Columns.Item(5).ColumnWidth = 21.035714285714285
$ws.Columns.Item(9).ColumnWidth = 23.660714285714285

# --- Selection moved from I14 to H12 ---------------------------------------
[void]$ws.Activate()
[void]$ws.Range("H12").Select()

Write-Output "done"
